$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.972010850906372
$ws.Range("B1").Value = 2.243095397949219
$ws.Range("C1").Value = 2.191790580749512
$ws.Range("D1").Value = 2.622453689575195
$ws.Range("E1").Value = 1.715325713157654
